# Insert a new data row before row 91 (new week's price record),
# shifting the existing rows 91-117 down to 92-118.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(91).Insert()

# Copy the row template (A,B,C,E,F,G,H,I,Q,R are identical across all
# data rows in this block) from the row that is now directly below
# (old row 91, now row 92) so formatting/styles line up, then overwrite
# the cells that hold this new week's figures.
$ws.Range("A92:R92").Copy() | Out-Null
$ws.Range("A91").PasteSpecial() | Out-Null
$excel.CutCopyMode = $false

$ws.Range("D91").Value = 44964
$ws.Range("J91").Value = 80
$ws.Range("K91").Value = 25000
$ws.Range("L91").Value = 25000
$ws.Range("M91").Value = 25000
$ws.Range("N91").Value = '$/saco 25 kilos'
$ws.Range("O91").Value = 'Región de La Araucanía'
$ws.Range("P91").Value = 1000
$ws.Range("Q91").Value = 25
